$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Last Name "
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Restults for URL1 "

# Update row 2
$ws.Range("A2").Value = "Moore"
$ws.Range("B2").Value = "Marie"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "headword"

# Add row 3
$ws.Range("C3").Value = "Yes"

# Set column C width (target stored OOXML width 26.5703125 characters;
# closest achievable value given pixel-rounding is 26.5)
$ws.Columns.Item(3).ColumnWidth = 25.67
